$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Add the two new Sterling Ratio test rows right after the existing data (row 82 -> 83, 84)
$ws.Cells.Item(83, 1).Value = "Sterling Ratio1"
$ws.Cells.Item(83, 2).Value = "Test Sterling ratio with scale=4"
$ws.Cells.Item(83, 3).Value = "Sterling_Ratio_test1"

$ws.Cells.Item(84, 1).Value = "Sterling Ratio2"
$ws.Cells.Item(84, 2).Value = "Test Sterling ratio with scale=12"
$ws.Cells.Item(84, 3).Value = "Sterling_Ratio_test2"

# Set the final selection to match the resulting cursor position
$ws.Range("J83").Select()
